# "vse je finally ready" - fix the character's name Gooper -> Cooper
# (two occurrences) and restore the missing space before the en-dash
# that follows the first occurrence ("Cooper a Mae" <space> "– ...").

$d = $word.ActiveDocument

# --- Occurrence 1: bold character-list entry "Gooper a Mae" ------------
# Replace the misspelled name; this run is immediately followed (with no
# run text change needed) by a separate run holding "– pokrytecti, ...".
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Gooper a Mae", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "Cooper a Mae", 1)

if ($found1) {
    # Insert the missing space between "Cooper a Mae" and the en-dash
    # run that follows it, keeping the bold formatting of the name run.
    $gap = $d.Range($rng1.End, $rng1.End)
    $gap.InsertAfter(" ")
    $gap.Font.Bold = $true
}

# --- Occurrence 2: plain-text synopsis "syn Gooper s manzelkou" --------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("syn Gooper s", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "syn Cooper s", 1)

Write-Host "Gooper->Cooper (char list): $found1"
Write-Host "Gooper->Cooper (synopsis): $found2"
